# Follow_Up_Log.xlsx maintenance edit
# - Fix a mis-entered FollowUpDate in row 8 (was stored as the text
#   "1/25/2016" carried over from a copy/paste, should be an actual date
#   serial like the surrounding rows) and give it the same date style.
# - Move the active cell/selection down to D11 (where work continued).
# - Nudge a few column widths back to their "best fit" size now that the
#   wide text string is gone from column D.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Fix the FollowUpDate value stored in D8 -------------------------------
# Previously a text string ("1/25/2016", style 1). Replace with the correct
# numeric date serial (11/3/2017) and apply the same date style used by the
# neighboring cells in column D (MM/DD/YY).
$ws.Range("D8").Value = 43042
$ws.Range("D8").NumberFormat = "MM/DD/YY"

# --- Move the selection/active cell ----------------------------------------
[void]$ws.Range("D11").Select()

# --- Re-tighten a few "best fit" column widths ------------------------------
# These columns use Calc/Excel's automatic (non-custom) width; after the fix
# above the optimum width shrank slightly for columns A-F.
$ws.Columns.Item(1).ColumnWidth = 8.666666666666666
$ws.Columns.Item(2).ColumnWidth = 12.333333333333334
$ws.Columns.Item(3).ColumnWidth = 12.666666666666666
$ws.Columns.Item(4).ColumnWidth = 12
$ws.Columns.Item(5).ColumnWidth = 12.333333333333334
$ws.Columns.Item(6).ColumnWidth = 16.166666666666668
